$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (Total) sheet: shift the existing data row down and insert a new
#    row for the 2022-Q3 summary.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Move the current row 2 (2021-Q3 totals) down to row 3, carrying formatting.
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("B2").Copy($total.Range("B3"))
$total.Range("C2").Copy($total.Range("C3"))
$total.Range("D2").Copy($total.Range("D3"))
$total.Range("A3").Value2 = 1

# Overwrite row 2 with the new 2022-Q3 summary values (A2 stays 0).
$total.Range("B2").Value2 = "2022-Q3"
$total.Range("C2").Value2 = 1
$total.Range("D2").Value2 = 2.39

# ---------------------------------------------------------------------------
# 2. Split the existing "2021-Q3" detail sheet into two sheets:
#      - a duplicate that keeps all the original 2021-Q3 fund data, and
#      - the original sheet object, which is repurposed (cleared + renamed)
#        to hold the brand-new 2022-Q3 fund data.
#    This keeps the original sheet's internal identity (and therefore its
#    position right after "总计") for the new "2022-Q3" data, while the
#    duplicate becomes the new "2021-Q3" sheet.
# ---------------------------------------------------------------------------
$origQ3 = $wb.Worksheets.Item(2)

$origQ3.Copy($null, $origQ3)
$dupQ3 = $wb.Worksheets.Item(3)

$origQ3.Name = "2022-Q3"
$dupQ3.Name = "2021-Q3"

$q2022 = $origQ3
$q2022.Cells.Clear()

# Match the page margins used on "总计" (points = inches * 72).
$q2022.PageSetup.LeftMargin = 54
$q2022.PageSetup.RightMargin = 54
$q2022.PageSetup.TopMargin = 72
$q2022.PageSetup.BottomMargin = 72
$q2022.PageSetup.HeaderMargin = 36
$q2022.PageSetup.FooterMargin = 36

# Header row (reuse the bold/bordered header style already used on "总计").
$total.Range("B1").Copy($q2022.Range("B1:H1"))
$q2022.Range("B1").Value2 = "基金代码"
$q2022.Range("C1").Value2 = "基金名称"
$q2022.Range("D1").Value2 = "基金规模"
$q2022.Range("E1").Value2 = "股票总仓位"
$q2022.Range("F1").Value2 = "仓位占比"
$q2022.Range("G1").Value2 = "持有市值(亿元)"
$q2022.Range("H1").Value2 = "仓位排名"

# Data row.
$total.Range("A2").Copy($q2022.Range("A2"))
$q2022.Range("A2").Value2 = 0
$q2022.Range("B2").Value2 = "'118001"
$q2022.Range("C2").Value2 = "易方达亚洲精选股票（QDII）"
$q2022.Range("D2").Value2 = "'46.17"
$q2022.Range("E2").Value2 = "'94.52"
$q2022.Range("F2").Value2 = "'5.17"
$q2022.Range("G2").Value2 = "'2.3870"
$q2022.Range("H2").Value2 = 10
